$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.807.09"
$ws.Range("E2").Value = "  +4.13%  "

$ws.Range("D3").Value = "3.623.88"
$ws.Range("E3").Value = "  +2.80%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "629.37"
$ws.Range("E5").Value = "  +3.31%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.16"
$ws.Range("E6").Value = "  +4.56%  "

$ws.Range("D7").Value = "3.622.34"
$ws.Range("E7").Value = "  +2.83%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.497"
$ws.Range("E9").Value = "  +2.83%  "

$ws.Range("E10").Value = "  +6.11%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.37"
$ws.Range("E11").Value = "  +6.59%  "

$ws.Range("E12").Value = "  +3.31%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000229"
$ws.Range("E13").Value = "  +4.04%  "

$ws.Range("E14").Value = "  +5.57%  "

$ws.Range("D15").Value = "4.236.93"
$ws.Range("E15").Value = "  +2.75%  "

$ws.Range("D16").Value = "3.629.13"
$ws.Range("E16").Value = "  +2.42%  "

$ws.Range("D17").Value = "69.877.30"
$ws.Range("E17").Value = "  +4.29%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.69"
$ws.Range("E19").Value = "  +5.72%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.02"
$ws.Range("E20").Value = "  +4.41%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.19"
$ws.Range("E21").Value = "  +13.21%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "463.86"
$ws.Range("E22").Value = "  +4.15%  "

$ws.Range("E23").Value = "  +2.53%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.97"
$ws.Range("E24").Value = "  +2.02%  "

$ws.Range("E25").Value = "  +11.27%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.74"
$ws.Range("E26").Value = "  +5.41%  "

$ws.Range("D27").Value = "3.769.05"
$ws.Range("E27").Value = "  +2.85%  "

$ws.Range("E28").Value = "  -0.05%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.28"
$ws.Range("E29").Value = "  +13.05%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.64"
$ws.Range("E30").Value = "  +4.63%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.72"
$ws.Range("E31").Value = "  +5.86%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.178"
$ws.Range("E32").Value = "  +11.64%  "

$ws.Range("E33").Value = "  +7.26%  "

$ws.Range("E34").Value = "  +0.08%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.97"
$ws.Range("E35").Value = "  +5.44%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "26.66"
$ws.Range("E36").Value = "  +3.74%  "

$ws.Range("D37").Value = "3.621.89"
$ws.Range("E37").Value = "  +3.05%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.50"
$ws.Range("E38").Value = "  +5.89%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.45"
$ws.Range("E39").Value = "  +13.46%  "

$ws.Range("E40").Value = "  -0.07%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0934"
$ws.Range("E41").Value = "  +8.05%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "179.81"
$ws.Range("E42").Value = "  +3.99%  "

$ws.Range("E43").Value = "  +0.07%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.70"
$ws.Range("E44").Value = "  +2.58%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "32.39"
$ws.Range("E45").Value = "  +18.82%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.917"
$ws.Range("E46").Value = "  +2.94%  "

$ws.Range("E47").Value = "  +11.92%  "

$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "46.32"
$ws.Range("E48").Value = "  +2.37%  "

$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.78"
$ws.Range("E49").Value = "  +9.40%  "

$ws.Range("E50").Value = "  +3.35%  "

$ws.Range("E51").Value = "  +8.76%  "
